# Weekly update for the Papaya price sheet (Vega Central Mapocho de Santiago).
# A new week of data (rows for 2023-07-18, serial 45125) is inserted above the
# existing history, pushing the previous rows down by four positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 8; Excel shifts rows 8:23 down to 12:27
# and carries their values/formatting (incl. the date style on column D) with
# them, exactly like the native "Insert Sheet Rows" command would.
$ws.Rows("8:11").Insert()

# Columns A,B,C,E,F,G,H,I,J,K,R are constant for every Papaya record in this
# sheet/origin combination.
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonE = 13
$commonF = "Fruta"
$commonG = 100108
$commonH = "Tropicales y subtropicales"
$commonI = 100108004
$commonJ = "Papaya"
$commonK = "Cultivar IV Región"
$commonR = "Provincia del Elquí"

function Set-PapayaRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $unidad, $precioKg, $kgUnidad)

    $ws.Cells.Item($row, 1).Value = $commonA
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    $ws.Cells.Item($row, 11).Value = $commonK
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $commonR
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-PapayaRow 8  45125 "Especial" 55 24000 24000 24000 "`$/bandeja 10 kilos" 2400 10
Set-PapayaRow 9  45125 "Primera"  60 20000 20000 20000 "`$/bandeja 10 kilos" 2000 10
Set-PapayaRow 10 45125 "Segunda"  60 15000 15000 15000 "`$/bandeja 10 kilos" 1500 10
Set-PapayaRow 11 45125 "Tercera"  45 12000 12000 12000 "`$/bandeja 10 kilos" 1200 10
